$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update date/time values (review meeting updates)
$ws.Range("A2").Value = 43803.447222222225
$ws.Range("B2").Value = 43803.470138888886
$ws.Range("C2").Value = 43805.599999999999

$ws.Range("A3").Value = 43802.470138888886
$ws.Range("B3").Value = 43802.568749999999
$ws.Range("C3").Value = 43803.5

$ws.Range("A4").Value = 43802.6
$ws.Range("B4").Value = 43803.367361111108
$ws.Range("C4").Value = 43804.916666666664

# Update the active cell selection from C4 to A4
$ws.Range("A4").Select()
